$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.903.02'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.814.15'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.13'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4663'
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3660'
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8688'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.30'
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").Value = '1.810.06'
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.387'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07081'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.516'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008706'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.65'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '26.928.27'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.296'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").Value = '2.036.72'
$ws.Range("E24").Value = '  -3.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.894'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.86'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.37'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.151'
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.37'
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08899'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7545'
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.156'
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.488'
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.912'
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05276'
$ws.Range("E38").Value = '  +0.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01947'
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.957'
$ws.Range("E40").Value = '  +0.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.225'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5303'
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.281'
$ws.Range("E43").Value = '  -2.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1650'
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.416'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("E46").Value = '  -2.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.36'
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.17'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.659'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06291'
$ws.Range("E51").Value = '  +0.07%  '
